$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of the adjacent header cell (G1) into the new H1 header
# cell so it gets the same bold/border/centered style, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save values for rows 2-14 (column H)
$saveValues = @(0, 1, 0, 0, 0, 1, 1, 1, 0, 1, 0, 1, 1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
